$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J4").ClearContents()
$ws.Range("J4").Borders.Item(7).LineStyle = -4142
$ws.Range("J4").Borders.Item(8).LineStyle = -4142
$ws.Range("J4").Borders.Item(9).LineStyle = -4142
$ws.Range("J4").Borders.Item(10).LineStyle = -4142
